# "Segunda optimización de Kp" - update the km input (C3) used by the Kp
# calculator and leave the selection where the author left it (K4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kp")
$ws.Activate()

# Second optimization pass: km goes from 6400.1 down to 1600, which ripples
# through the Mp/Kp/kp*75 columns (I3:J11) via the existing formulas.
$ws.Range("C3").Value = 1600

# Leave the selection on K4, matching the author's saved view.
$ws.Range("K4").Select()
